$d = $word.ActiveDocument

$replacements = @(
    @{old="73×41=2993"; new="49×58=2842"},
    @{old="39×71=2769"; new="46×70=3220"},
    @{old="11×28=308";  new="77×70=5390"},
    @{old="12×83=996";  new="65×16=1040"},
    @{old="16×13=208";  new="72×62=4464"},
    @{old="50×25=1250"; new="61×51=3111"},
    @{old="62×51=3162"; new="95×99=9405"},
    @{old="90×65=5850"; new="73×71=5183"},
    @{old="55×41=2255"; new="12×71=852"},
    @{old="30×38=1140"; new="55×22=1210"},
    @{old="13×48=624";  new="24×98=2352"},
    @{old="65×61=3965"; new="21×42=882"},
    @{old="15×78=1170"; new="66×69=4554"},
    @{old="56×85=4760"; new="40×68=2720"},
    @{old="31×61=1891"; new="43×42=1806"},
    @{old="20×77=1540"; new="50×19=950"},
    @{old="36×82=2952"; new="74×62=4588"},
    @{old="67×47=3149"; new="55×85=4675"},
    @{old="57×85=4845"; new="93×66=6138"},
    @{old="89×17=1513"; new="30×81=2430"},
    @{old="99×59=5841"; new="27×69=1863"},
    @{old="13×59=767";  new="29×96=2784"},
    @{old="64×17=1088"; new="22×14=308"},
    @{old="63×85=5355"; new="62×12=744"},
    @{old="92×89=8188"; new="41×39=1599"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
